$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = '24.863.31'
$cell.ClearFormats()
$ws.Cells.Item(2, 5).Value = '  +0.28%  '

# Row 3
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.706.12'
$cell.ClearFormats()
$ws.Cells.Item(3, 5).Value = '  +0.20%  '

# Row 4
$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.000'
$cell.ClearFormats()
$ws.Cells.Item(4, 5).Value = '  -0.55%  '

# Row 5
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '318.14'
$cell.ClearFormats()
$ws.Cells.Item(5, 5).Value = '  +0.18%  '

# Row 6
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.ClearFormats()

# Row 7
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.3973'
$cell.ClearFormats()
$ws.Cells.Item(7, 5).Value = '  +0.27%  '

# Row 8
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.4074'
$cell.ClearFormats()
$ws.Cells.Item(8, 5).Value = '  -0.48%  '

# Row 9
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.491'
$cell.ClearFormats()
$ws.Cells.Item(9, 5).Value = '  -1.02%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  -0.70%  '

# Row 11
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '53.62'
$cell.ClearFormats()
$ws.Cells.Item(11, 5).Value = '  +1.38%  '

# Row 12
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.08830'
$cell.ClearFormats()
$ws.Cells.Item(12, 5).Value = '  -1.25%  '

# Row 13
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '26.42'
$cell.ClearFormats()
$ws.Cells.Item(13, 5).Value = '  +7.97%  '

# Row 14
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.526'
$cell.ClearFormats()
$ws.Cells.Item(14, 5).Value = '  -2.29%  '

# Row 15
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.186'
$cell.ClearFormats()
$ws.Cells.Item(15, 5).Value = '  +0.28%  '

# Row 16
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.00001360'
$cell.ClearFormats()
$ws.Cells.Item(16, 5).Value = '  +1.79%  '

# Row 17
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.665.30'
$cell.ClearFormats()
$ws.Cells.Item(17, 5).Value = '  -2.52%  '

# Row 18
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '96.91'
$cell.ClearFormats()
$ws.Cells.Item(18, 5).Value = '  -3.27%  '

# Row 19
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.07191'
$cell.ClearFormats()
$ws.Cells.Item(19, 5).Value = '  +0.69%  '

# Row 20
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '21.13'
$cell.ClearFormats()
$ws.Cells.Item(20, 5).Value = '  +5.22%  '

# Row 21
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.301'
$cell.ClearFormats()
$ws.Cells.Item(21, 5).Value = '  -0.24%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  -0.49%  '

# Row 23
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '14.43'
$cell.ClearFormats()
$ws.Cells.Item(23, 5).Value = '  -0.90%  '

# Row 24
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '24.861.97'
$cell.ClearFormats()
$ws.Cells.Item(24, 5).Value = '  +0.31%  '

# Row 25
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.992'
$cell.ClearFormats()
$ws.Cells.Item(25, 5).Value = '  -3.15%  '

# Row 26
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.331'
$cell.ClearFormats()
$ws.Cells.Item(26, 5).Value = '  -0.40%  '

# Row 27
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '23.27'
$cell.ClearFormats()
$ws.Cells.Item(27, 5).Value = '  +0.87%  '

# Row 28
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.265'
$cell.ClearFormats()
$ws.Cells.Item(28, 5).Value = '  +20.28%  '

# Row 29
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '167.02'
$cell.ClearFormats()
$ws.Cells.Item(29, 5).Value = '  +1.01%  '

# Row 30
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '145.88'
$cell.ClearFormats()
$ws.Cells.Item(30, 5).Value = '  +4.47%  '

# Row 31
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.395'
$cell.ClearFormats()
$ws.Cells.Item(31, 5).Value = '  -10.98%  '

# Row 32
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.242'
$cell.ClearFormats()
$ws.Cells.Item(32, 5).Value = '  +13.82%  '

# Row 33
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.860.94'
$cell.ClearFormats()
$ws.Cells.Item(33, 5).Value = '  -1.88%  '

# Row 34
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.08813'
$cell.ClearFormats()
$ws.Cells.Item(34, 5).Value = '  -4.08%  '

# Row 35
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.03221'
$cell.ClearFormats()
$ws.Cells.Item(35, 5).Value = '  +4.83%  '

# Row 36
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.214'
$cell.ClearFormats()
$ws.Cells.Item(36, 5).Value = '  -10.80%  '

# Row 37
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.040'
$cell.ClearFormats()
$ws.Cells.Item(37, 5).Value = '  -4.19%  '

# Row 38
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.2899'
$cell.ClearFormats()
$ws.Cells.Item(38, 5).Value = '  +2.88%  '

# Row 39
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.8568'
$cell.ClearFormats()
$ws.Cells.Item(39, 5).Value = '  +9.22%  '

# Row 40
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '10.93'
$cell.ClearFormats()
$ws.Cells.Item(40, 5).Value = '  -1.61%  '

# Row 41
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.09295'
$cell.ClearFormats()
$ws.Cells.Item(41, 5).Value = '  -0.23%  '

# Row 42
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '14.16'
$cell.ClearFormats()
$ws.Cells.Item(42, 5).Value = '  -2.98%  '

# Row 43
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.474'
$cell.ClearFormats()
$ws.Cells.Item(43, 5).Value = '  -0.47%  '

# Row 44
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '17.55'
$cell.ClearFormats()
$ws.Cells.Item(44, 5).Value = '  +6.93%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'NEARProtocol'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.715'
$cell.ClearFormats()
$ws.Cells.Item(45, 5).Value = '  +2.81%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'Decentraland'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.7484'
$cell.ClearFormats()
$ws.Cells.Item(46, 5).Value = '  +2.27%  '

# Row 47
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.265'
$cell.ClearFormats()
$ws.Cells.Item(47, 5).Value = '  +0.24%  '

# Row 48
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.400'
$cell.ClearFormats()
$ws.Cells.Item(48, 5).Value = '  +3.11%  '

# Row 49
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.9993'
$cell.ClearFormats()
$ws.Cells.Item(49, 5).Value = '  -0.43%  '

# Row 50
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '141.31'
$cell.ClearFormats()
$ws.Cells.Item(50, 5).Value = '  -0.07%  '

# Row 51
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.08343'
$cell.ClearFormats()
$ws.Cells.Item(51, 5).Value = '  +3.38%  '
